$d = $word.ActiveDocument

# Paragraph 5: "Tìm các biểu mẫu và quy định: Trương Việt Hoàng"
#           -> "Tìm các biểu mẫu và quy định, trình bày các công nghệ: Trương Việt Hoàng"
$d.Content.Find.Execute(
    "Tìm các biểu mẫu và quy định: Trương Việt Hoàng",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Tìm các biểu mẫu và quy định, trình bày các công nghệ: Trương Việt Hoàng",
    2
)

# Paragraph 6: "BFD và DFD (mức 0, 1): Đoàn Duy Khánh"
#           -> "BFD và DFD (mức 0, 1), usecase diagram, class diagram: Đoàn Duy Khánh"
$d.Content.Find.Execute(
    "BFD và DFD (mức 0, 1): Đoàn Duy Khánh",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "BFD và DFD (mức 0, 1), usecase diagram, class diagram: Đoàn Duy Khánh",
    2
)
